# Auto-generated edit script: applies numeric value updates to the
# "Adamantoise_Profits" market-data workbook across all 8 job sheets,
# matching the scheduled-runner commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1425.3549
$ws.Range("I15").Value = 1425.3549
$ws.Range("K15").Value = 4276.0647
$ws.Range("M15").Value = -4107.0647
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("H51").Value = 3823.1155
$ws.Range("J51").Value = 2960.1333
$ws.Range("L51").Value = 2960.1333
$ws.Range("N51").Value = -3928.1333
$ws.Range("H96").Value = 2185.2632
$ws.Range("J96").Value = 2904.2
$ws.Range("L96").Value = 8712.599999999999
$ws.Range("N96").Value = -11458.6
$ws.Range("H97").Value = 1724.6428
$ws.Range("J97").Value = 1724.6428
$ws.Range("L97").Value = 5173.928400000001
$ws.Range("N97").Value = -6165.928400000001
$ws.Range("H98").Value = 953.05457
$ws.Range("I98").Value = 953.05457
$ws.Range("K98").Value = 953.05457
$ws.Range("M98").Value = 544.94543
$ws.Range("H100").Value = 2610.5557
$ws.Range("I100").Value = 1438.4
$ws.Range("J100").Value = 4075.75
$ws.Range("K100").Value = 1438.4
$ws.Range("L100").Value = 4075.75
$ws.Range("M100").Value = -897.4000000000001
$ws.Range("N100").Value = -5157.75
$ws.Range("H101").Value = 397
$ws.Range("I101").Value = 397
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1191
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 431
$ws.Range("H111").Value = 1991.6666
$ws.Range("I111").Value = 1991.6666
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 5974.9998
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -2907.9998
$ws.Range("H118").Value = 1733.3
$ws.Range("I118").Value = 1010.5714
$ws.Range("J118").Value = 3419.6667
$ws.Range("K118").Value = 3031.7142
$ws.Range("L118").Value = 10259.0001
$ws.Range("M118").Value = -1374.7142
$ws.Range("N118").Value = -13573.0001
$ws.Range("H122").Value = 953.05457
$ws.Range("I122").Value = 953.05457
$ws.Range("K122").Value = 2859.16371
$ws.Range("M122").Value = -409.1637099999998
$ws.Range("H129").Value = 1621.5883
$ws.Range("I129").Value = 670.875
$ws.Range("K129").Value = 2012.625
$ws.Range("M129").Value = 2987.375
$ws.Range("N43").ClearContents()
$ws.Range("N101").ClearContents()
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10968392
$ws.Range("I32").Value = 6411803.5
$ws.Range("K32").Value = 6411803.5
$ws.Range("M32").Value = -6411516.5
$ws.Range("H80").Value = 200000
$ws.Range("J80").Value = 200000
$ws.Range("L80").Value = 200000
$ws.Range("N80").Value = -201996
$ws.Range("H83").Value = 200000
$ws.Range("J83").Value = 200000
$ws.Range("L83").Value = 600000
$ws.Range("N83").Value = -609984
$ws.Range("H97").Value = 816.3871
$ws.Range("I97").Value = 625.38464
$ws.Range("K97").Value = 625.38464
$ws.Range("M97").Value = -129.38464
$ws.Range("H122").Value = 1470.56
$ws.Range("I122").Value = 1174.2391
$ws.Range("K122").Value = 3522.7173
$ws.Range("M122").Value = -1072.7173

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1669.6
$ws.Range("I86").Value = 2767.6667
$ws.Range("J86").Value = 1199
$ws.Range("K86").Value = 2767.6667
$ws.Range("L86").Value = 1199
$ws.Range("M86").Value = -1644.6667
$ws.Range("N86").Value = -3445
$ws.Range("H89").Value = 1669.6
$ws.Range("I89").Value = 2767.6667
$ws.Range("J89").Value = 1199
$ws.Range("K89").Value = 13838.3335
$ws.Range("L89").Value = 5995
$ws.Range("M89").Value = -8222.333500000001
$ws.Range("N89").Value = -17227
$ws.Range("H93").Value = 51505
$ws.Range("J93").Value = 51505
$ws.Range("L93").Value = 51505
$ws.Range("N93").Value = -55249
$ws.Range("H94").Value = 969.2083
$ws.Range("I94").Value = 479.7647
$ws.Range("K94").Value = 479.7647
$ws.Range("M94").Value = -28.7647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3569
$ws.Range("I2").Value = 3330.5
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 3330.5
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -3217.5
$ws.Range("N2").Value = -5226
$ws.Range("H6").Value = 4338282.5
$ws.Range("I6").Value = 4338282.5
$ws.Range("K6").Value = 4338282.5
$ws.Range("M6").Value = -4338169.5
$ws.Range("H16").Value = 1807.375
$ws.Range("I16").Value = 1491.8
$ws.Range("K16").Value = 1491.8
$ws.Range("M16").Value = -1204.8
$ws.Range("H31").Value = 3870.6206
$ws.Range("I31").Value = 2094.3572
$ws.Range("J31").Value = 5528.467
$ws.Range("K31").Value = 2094.3572
$ws.Range("L31").Value = 5528.467
$ws.Range("M31").Value = -1799.3572
$ws.Range("N31").Value = -6118.467
$ws.Range("H34").Value = 3870.6206
$ws.Range("I34").Value = 2094.3572
$ws.Range("J34").Value = 5528.467
$ws.Range("K34").Value = 2094.3572
$ws.Range("L34").Value = 5528.467
$ws.Range("M34").Value = -1892.3572
$ws.Range("N34").Value = -5932.467
$ws.Range("H103").Value = 60769.363
$ws.Range("I103").Value = 14624
$ws.Range("J103").Value = 87138.14
$ws.Range("K103").Value = 14624
$ws.Range("L103").Value = 87138.14
$ws.Range("M103").Value = -13452
$ws.Range("N103").Value = -89482.14
$ws.Range("H113").Value = 1807.375
$ws.Range("I113").Value = 1491.8
$ws.Range("K113").Value = 1491.8
$ws.Range("M113").Value = 678.2
$ws.Range("H122").Value = 2578.2856
$ws.Range("I122").Value = 1512.25
$ws.Range("J122").Value = 3999.6667
$ws.Range("K122").Value = 4536.75
$ws.Range("L122").Value = 11999.0001
$ws.Range("M122").Value = -2086.75
$ws.Range("N122").Value = -16899.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("H64").Value = 875
$ws.Range("I64").Value = 875
$ws.Range("K64").Value = 2625
$ws.Range("M64").Value = -2355
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("H67").Value = 875
$ws.Range("I67").Value = 875
$ws.Range("K67").Value = 2625
$ws.Range("M67").Value = -1689
$ws.Range("H132").Value = 1713.5385
$ws.Range("J132").Value = 1713.5385
$ws.Range("L132").Value = 15421.8465
$ws.Range("N132").Value = -20481.8465
$ws.Range("H134").Value = 3011
$ws.Range("I134").Value = 3011
$ws.Range("K134").Value = 9033
$ws.Range("M134").Value = -3963
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1501.4897
$ws.Range("I102").Value = 1435
$ws.Range("K102").Value = 1435
$ws.Range("M102").Value = 187

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13882.828
$ws.Range("I7").Value = 13684
$ws.Range("J7").Value = 14379.9
$ws.Range("K7").Value = 13684
$ws.Range("L7").Value = 14379.9
$ws.Range("M7").Value = -13572
$ws.Range("N7").Value = -14603.9
$ws.Range("H40").Value = 6974.4136
$ws.Range("I40").Value = 7094.2085
$ws.Range("J40").Value = 6399.4
$ws.Range("K40").Value = 7094.2085
$ws.Range("L40").Value = 6399.4
$ws.Range("M40").Value = -6958.2085
$ws.Range("N40").Value = -6671.4
$ws.Range("H93").Value = 43479188
$ws.Range("I93").Value = 58824236
$ws.Range("K93").Value = 58824236
$ws.Range("M93").Value = -58822988
$ws.Range("H122").Value = 3580.8125
$ws.Range("I122").Value = 3580.8125
$ws.Range("K122").Value = 10742.4375
$ws.Range("M122").Value = -8292.4375
$ws.Range("H126").Value = 13882.828
$ws.Range("I126").Value = 13684
$ws.Range("J126").Value = 14379.9
$ws.Range("K126").Value = 41052
$ws.Range("L126").Value = 43139.7
$ws.Range("M126").Value = -38582
$ws.Range("N126").Value = -48079.7
$ws.Range("H133").Value = 49888
$ws.Range("J133").Value = 49888
$ws.Range("L133").Value = 49888
$ws.Range("N133").Value = -54948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1463.5294
$ws.Range("I100").Value = 1754.1818
$ws.Range("K100").Value = 3508.3636
$ws.Range("M100").Value = -2967.3636
$ws.Range("H122").Value = 2781.3157
$ws.Range("I122").Value = 2077
$ws.Range("K122").Value = 6231
$ws.Range("M122").Value = -3781
$ws.Range("H126").Value = 4788.1177
$ws.Range("I126").Value = 5031.125
$ws.Range("K126").Value = 15093.375
$ws.Range("M126").Value = -12623.375
$ws.Range("H132").Value = 5022.95
$ws.Range("I132").Value = 4199.6665
$ws.Range("J132").Value = 5696.5454
$ws.Range("K132").Value = 12598.9995
$ws.Range("L132").Value = 17089.6362
$ws.Range("M132").Value = -10068.9995
$ws.Range("N132").Value = -22149.6362
$ws.Range("H133").Value = 79999
$ws.Range("J133").Value = 79999
$ws.Range("L133").Value = 79999
$ws.Range("N133").Value = -90119
